# Update ObjTables header strings: bump the timestamp and add tableFormat='row'
# to every per-table "!!ObjTables type='Data' ..." string, and bump the
# timestamp on the workbook-level "!!!ObjTables ..." string (sheet 1, A1).

$wb = $excel.ActiveWorkbook

$oldDate = "2020-03-09 13:01:23"
$newDate = "2020-03-09 15:31:44"

$ids = @(
    "Compartment",
    "Compound",
    "Definition",
    "Enzyme",
    "FbcObjective",
    "Gene",
    "Layout",
    "Measurement",
    "PbConfig",
    "Position",
    "Protein",
    "Quantity",
    "QuantityInfo",
    "QuantityMatrix",
    "Reaction",
    "ReactionStoichiometry",
    "Regulator",
    "Relation",
    "Relationship",
    "SparseMatrix",
    "SparseMatrixColumn",
    "SparseMatrixOrdered",
    "SparseMatrixRow",
    "StoichiometricMatrix",
    "rxnconContingencyList",
    "rxnconReactionList"
)

for ($i = 0; $i -lt $ids.Count; $i++) {
    $sheetIndex = $i + 1
    $id = $ids[$i]
    $ws = $wb.Worksheets.Item($sheetIndex)

    $wasProtected = $ws.ProtectContents
    if ($wasProtected) {
        $ws.Unprotect()
    }

    if ($sheetIndex -eq 1) {
        # First sheet also carries the workbook-level "!!!ObjTables ..." line in A1.
        $topCell = $ws.Range("A1")
        $topCell.Value2 = $topCell.Value2.Replace($oldDate, $newDate)

        $headerCell = $ws.Range("A2")
    } else {
        $headerCell = $ws.Range("A1")
    }

    $newValue = "!!ObjTables type='Data' id='$id' name='$id' date='$newDate' objTablesVersion='0.0.8' tableFormat='row'"
    $headerCell.Value2 = $newValue

    if ($wasProtected) {
        $ws.Protect($null, $true, $true, $true)
    }
}
